# Apply the data refresh: new header wording, reshuffled/updated rows 2-7,
# and a brand-new row 8 (sheet grows from A1:G7 to A1:G8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Vorname"
$ws.Range("B1").Value = "Nachname"

# --- Data rows (Vorname, Nachname, Klasse, Ankunft, Verlassen, DauerMinuten) ---
$data = @(
    @("Stephan", "Fuchs",             "2020B", "27.10.2024 18:17", "27.10.2024 20:28", 131),
    @("Steven",  "Mustermann",        "2013A", "27.10.2024 18:18", "27.10.2024 18:24", 6),
    @("Maike",   "perfect",           "2010B", "27.10.2024 18:19", "27.10.2024 18:24", 5),
    @("Max",     "Schmitz",           "2020A", "27.10.2024 18:20", "27.10.2024 18:20", 0),
    @("Otto",    "Langnamenokidoki",  "2015B", "27.10.2024 18:52", "27.10.2024 19:52", 60),
    @("Maike",   "perfect",           "2010B", "14.03.2025 18:35", "14.03.2025 21:36", 181),
    @("Max",     "Testico",           "2015A", "14.03.2025 19:05", "14.03.2025 20:03", 58)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}

# Row 8 is brand new; give it the same empty "Kommentar" cell the other
# data rows already carry in column G so the used range becomes A1:G8.
$ws.Cells.Item(8, 7).Value = ""
